$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Append a new row to the weekly-progress table (mirrors the author typing
#    a new entry for "6 / 7.6.2019 / 1.Modified Binary Tree Structure /
#    2. Structure Review / Done" under the existing last row).
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$rowIdx = $t.Rows.Count

# Column 1: "No"
$t.Cell($rowIdx, 1).Range.Text = "6"

# Column 2: "Date"
$t.Cell($rowIdx, 2).Range.Text = "7.6.2019"

# Column 3: "Assigned Tasks/ Description of work" -- two lines of text followed
# by a trailing blank paragraph (as in the existing rows of the table).
$descCell = $t.Cell($rowIdx, 3)
$descRng = $descCell.Range
$descRng.End = $descRng.End - 1
$descRng.Text = "1.Modified Binary Tree Structure" + [char]13 + "2. Structure Review" + [char]13

# Column 4: "Completion status"
$t.Cell($rowIdx, 4).Range.Text = "Done"

# Column 5: "Remarks" stays blank (left untouched, as cloned by Rows.Add()).

# ---------------------------------------------------------------------------
# 2. The "Supervisor's Comment" heading run carried a stale
#    <w:lastRenderedPageBreak/> marker left over from the previous
#    pagination; touching/re-writing the run drops it, same as Word does
#    when the surrounding content is re-flowed because of the edit above.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Supervisor" + [char]8217 + "s Comment ", $false, $false, $false, $false, $false, $true, 1, $false, "Supervisor" + [char]8217 + "s Comment ", 2) | Out-Null
